$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-08-12 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-08-13 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("918÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "100÷3=", 2) | Out-Null
$d.Content.Find.Execute("681÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "588÷2=", 2) | Out-Null
$d.Content.Find.Execute("148÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "171÷3=", 2) | Out-Null
$d.Content.Find.Execute("682÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "118÷9=", 2) | Out-Null
$d.Content.Find.Execute("430÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "107÷7=", 2) | Out-Null
$d.Content.Find.Execute("314÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "781÷5=", 2) | Out-Null
$d.Content.Find.Execute("335÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "219÷5=", 2) | Out-Null
$d.Content.Find.Execute("804÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "232÷9=", 2) | Out-Null
$d.Content.Find.Execute("935÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "567÷3=", 2) | Out-Null
$d.Content.Find.Execute("309÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "308÷4=", 2) | Out-Null
$d.Content.Find.Execute("328÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "779÷9=", 2) | Out-Null
$d.Content.Find.Execute("532÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "653÷7=", 2) | Out-Null
$d.Content.Find.Execute("364÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "446÷2=", 2) | Out-Null
$d.Content.Find.Execute("688÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "731÷6=", 2) | Out-Null
$d.Content.Find.Execute("436÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "614÷9=", 2) | Out-Null
$d.Content.Find.Execute("134÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "586÷8=", 2) | Out-Null
$d.Content.Find.Execute("595÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "174÷8=", 2) | Out-Null
$d.Content.Find.Execute("379÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "661÷5=", 2) | Out-Null
$d.Content.Find.Execute("149÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "759÷4=", 2) | Out-Null
$d.Content.Find.Execute("482÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "507÷8=", 2) | Out-Null
$d.Content.Find.Execute("549÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "308÷2=", 2) | Out-Null
$d.Content.Find.Execute("158÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "572÷7=", 2) | Out-Null
$d.Content.Find.Execute("333÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "759÷6=", 2) | Out-Null
$d.Content.Find.Execute("451÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "605÷9=", 2) | Out-Null
$d.Content.Find.Execute("380÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "120÷7=", 2) | Out-Null
